$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 update: "Store 1" -> "store 1", "S/MAX 22" GUARD" -> "Blade", D2 2000 -> 1626
$ws.Range("A2").Value = "store 1"
$ws.Range("B2").Value = "Blade"
$ws.Range("C2").Value = 2000
$ws.Range("D2").Value = 1626
$ws.Range("E2").Value = 1626

# New rows 3-11
$ws.Range("A3").Value = "SKICROP TRINITY PVT LTD"
$ws.Range("B3").Value = "Blade"
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 500
$ws.Range("E3").Value = 500

$ws.Range("A4").Value = "store 1"
$ws.Range("B4").Value = "Fitting Screw"
$ws.Range("C4").Value = 5000
$ws.Range("D4").Value = 5003
$ws.Range("E4").Value = 5003

$ws.Range("A5").Value = "SKICROP TRINITY PVT LTD"
$ws.Range("B5").Value = "Fitting Screw"
$ws.Range("C5").Value = 200
$ws.Range("D5").Value = 300
$ws.Range("E5").Value = 300

$ws.Range("A6").Value = "store 1"
$ws.Range("B6").Value = "PVC Socket"
$ws.Range("C6").Value = 5000
$ws.Range("D6").Value = 4800
$ws.Range("E6").Value = 4800

$ws.Range("A7").Value = "store 1"
$ws.Range("B7").Value = "Side Knob (S/Max, Max, Farata)"
$ws.Range("C7").Value = 2000
$ws.Range("D7").Value = 1700
$ws.Range("E7").Value = 1700

$ws.Range("A8").Value = "SKICROP TRINITY PVT LTD"
$ws.Range("B8").Value = "PVC Socket"
$ws.Range("C8").Value = 200
$ws.Range("D8").Value = 200
$ws.Range("E8").Value = 200

$ws.Range("A9").Value = "SKICROP TRINITY PVT LTD"
$ws.Range("B9").Value = "Side Knob (S/Max, Max, Farata)"
$ws.Range("C9").Value = 100
$ws.Range("D9").Value = 300
$ws.Range("E9").Value = 300

$ws.Range("A10").Value = "store 1"
$ws.Range("B10").Value = "Box"
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

$ws.Range("A11").Value = "store 1"
$ws.Range("B11").Value = "Pink Tape 6mm X 55 Mts"
$ws.Range("C11").Value = 33
$ws.Range("D11").Value = 33
$ws.Range("E11").Value = 33
